$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.919.93"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "2.298.48"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D4").Value = "`'0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "`'514.98"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "`'129.93"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("D7").Value = "`'0.998"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "`'0.526"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "2.304.60"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "`'0.0986"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "`'5.20"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").Value = "`'0.333"
$ws.Range("E13").Value = "  -2.39%  "
$ws.Range("D14").Value = "2.705.95"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "`'23.07"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").Value = "55.852.24"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "`'0.0000131"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "2.309.07"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "`'10.27"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("D20").Value = "`'325.44"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "`'4.09"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "`'6.64"
$ws.Range("E22").Value = "  +0.99%  "
$ws.Range("D23").Value = "`'0.997"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "`'60.60"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "`'8.46"
$ws.Range("E26").Value = "  +6.15%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "`'0.997"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "`'1.31"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").Value = "`'167.07"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "`'1.68"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "0.0₃0704"
$ws.Range("E31").Value = "  -4.63%  "
$ws.Range("D32").Value = "`'6.01"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "`'18.09"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").Value = "`'0.995"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "`'1.22"
$ws.Range("E36").Value = "  -2.55%  "
$ws.Range("D37").Value = "`'0.872"
$ws.Range("E37").Value = "  -4.95%  "
$ws.Range("D38").Value = "`'3.84"
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "`'1.55"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "`'38.27"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "`'146.87"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("D42").Value = "`'0.370"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "`'3.52"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").Value = "`'278.16"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "`'4.89"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("D46").Value = "`'0.0920"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").Value = "`'0.0491"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").Value = "`'0.551"
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("D49").Value = "`'17.76"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("D50").Value = "`'0.375"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "`'0.0211"
$ws.Range("E51").Value = "  -3.10%  "
